# Update NATMI Tnfsf14-Tnfrsf14 LR-pair sheet with recomputed TPM-based values.
# Ligand-expressing cell counts / detection rates / expression values and the
# resulting specificity & edge-weight metrics (columns E-T) were recalculated
# with the new TPM normalization; sending/receptor cluster labels (A-D) and
# receptor-expressing cell counts (K-L) are unchanged.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$updates = @{
    "E2" = 3
    "F2" = 1
    "G2" = 5.57196
    "H2" = 16.71588
    "I2" = 0.2041274261050298
    "J2" = 0.2041274261050299
    "M2" = 2.862697666666667
    "N2" = 8.588093000000001
    "O2" = 0.06997451156315468
    "P2" = 0.06997451156315466
    "Q2" = 15.95083689076
    "R2" = 143.55753201684
    "S2" = 0.01428371693834341
    "T2" = 0.01428371693834341
    "E3" = 3
    "F3" = 1
    "G3" = 5.57196
    "H3" = 16.71588
    "I3" = 0.2041274261050298
    "J3" = 0.2041274261050299
    "O3" = 0.06164554086142937
    "P3" = 0.06164554086142937
    "Q3" = 14.05223052448
    "R3" = 126.47007472032
    "S3" = 0.01258354558689602
    "T3" = 0.01258354558689602
    "E4" = 3
    "F4" = 1
    "G4" = 5.57196
    "H4" = 16.71588
    "I4" = 0.2041274261050298
    "J4" = 0.2041274261050299
    "M4" = 17.93719066666667
    "N4" = 53.811572
    "O4" = 0.4384487297873381
    "P4" = 0.4384487297873381
    "Q4" = 99.94530890703999
    "R4" = 899.5077801633599
    "S4" = 0.08949941069050905
    "T4" = 0.08949941069050907
    "E5" = 3
    "F5" = 1
    "G5" = 5.57196
    "H5" = 16.71588
    "I5" = 0.2041274261050298
    "J5" = 0.2041274261050299
    "M5" = 0.5052413333333333
    "N5" = 1.515724
    "O5" = 0.01234989497255689
    "P5" = 0.01234989497255689
    "Q5" = 2.81518449968
    "R5" = 25.33666049712
    "S5" = 0.002520952273415487
    "T5" = 0.002520952273415487
    "E6" = 3
    "F6" = 1
    "G6" = 5.57196
    "H6" = 16.71588
    "I6" = 0.2041274261050298
    "J6" = 0.2041274261050299
    "M6" = 11.46238066666667
    "N6" = 34.387142
    "O6" = 0.280181347070047
    "P6" = 0.280181347070047
    "Q6" = 63.86792657944
    "R6" = 574.81133921496
    "S6" = 0.05719269722004875
    "T6" = 0.05719269722004874
    "E7" = 3
    "F7" = 1
    "G7" = 5.57196
    "H7" = 16.71588
    "I7" = 0.2041274261050298
    "J7" = 0.2041274261050299
    "M7" = 5.621112333333333
    "N7" = 16.863337
    "O7" = 0.1373999757454739
    "P7" = 0.1373999757454739
    "Q7" = 31.32061307684
    "R7" = 281.88551769156
    "S7" = 0.02804710339581711
    "T7" = 0.02804710339581712
    "G8" = 20.56891366666667
    "H8" = 61.70674100000001
    "I8" = 0.7535372480335895
    "J8" = 0.7535372480335896
    "M8" = 2.862697666666667
    "N8" = 8.588093000000001
    "O8" = 0.06997451156315468
    "P8" = 0.06997451156315466
    "Q8" = 58.8825811594348
    "R8" = 529.9432304349131
    "S8" = 0.05272840087579416
    "T8" = 0.05272840087579416
    "G9" = 20.56891366666667
    "H9" = 61.70674100000001
    "I9" = 0.7535372480335895
    "J9" = 0.7535372480335896
    "O9" = 0.06164554086142937
    "P9" = 0.06164554086142937
    "Q9" = 51.87386780991378
    "R9" = 466.864810289224
    "S9" = 0.04645221121426368
    "T9" = 0.04645221121426368
    "G10" = 20.56891366666667
    "H10" = 61.70674100000001
    "I10" = 0.7535372480335895
    "J10" = 0.7535372480335896
    "M10" = 17.93719066666667
    "N10" = 53.811572
    "O10" = 0.4384487297873381
    "P10" = 0.4384487297873381
    "Q10" = 368.9485262452058
    "R10" = 3320.536736206852
    "S10" = 0.3303874492477736
    "T10" = 0.3303874492477736
    "G11" = 20.56891366666667
    "H11" = 61.70674100000001
    "I11" = 0.7535372480335895
    "J11" = 0.7535372480335896
    "M11" = 0.5052413333333333
    "N11" = 1.515724
    "O11" = 0.01234989497255689
    "P11" = 0.01234989497255689
    "Q11" = 10.39226536616489
    "R11" = 93.53038829548402
    "S11" = 0.009306105871124383
    "T11" = 0.009306105871124385
    "G12" = 20.56891366666667
    "H12" = 61.70674100000001
    "I12" = 0.7535372480335895
    "J12" = 0.7535372480335896
    "M12" = 11.46238066666667
    "N12" = 34.387142
    "O12" = 0.280181347070047
    "P12" = 0.280181347070047
    "Q12" = 235.7687183471358
    "R12" = 2121.918465124222
    "S12" = 0.2111270812215073
    "T12" = 0.2111270812215073
    "G13" = 20.56891366666667
    "H13" = 61.70674100000001
    "I13" = 0.7535372480335895
    "J13" = 0.7535372480335896
    "M13" = 5.621112333333333
    "N13" = 16.863337
    "O13" = 0.1373999757454739
    "P13" = 0.1373999757454739
    "Q13" = 115.6201742949686
    "R13" = 1040.581568654717
    "S13" = 0.1035359996031263
    "T13" = 0.1035359996031263
    "G14" = 1.155605333333333
    "H14" = 3.466816
    "I14" = 0.04233532586138062
    "J14" = 0.04233532586138063
    "M14" = 2.862697666666667
    "N14" = 8.588093000000001
    "O14" = 0.06997451156315468
    "P14" = 0.06997451156315466
    "Q14" = 3.308148691320889
    "R14" = 29.773338221888
    "S14" = 0.002962393749017099
    "T14" = 0.002962393749017099
    "G15" = 1.155605333333333
    "H15" = 3.466816
    "I15" = 0.04233532586138062
    "J15" = 0.04233532586138063
    "O15" = 0.06164554086142937
    "P15" = 0.06164554086142937
    "Q15" = 2.914384263224888
    "R15" = 26.22945836902399
    "S15" = 0.002609784060269666
    "T15" = 0.002609784060269667
    "G16" = 1.155605333333333
    "H16" = 3.466816
    "I16" = 0.04233532586138062
    "J16" = 0.04233532586138063
    "M16" = 17.93719066666667
    "N16" = 53.811572
    "O16" = 0.4384487297873381
    "P16" = 0.4384487297873381
    "Q16" = 20.72831319941688
    "R16" = 186.554818794752
    "S16" = 0.01856186984905538
    "T16" = 0.01856186984905538
    "G17" = 1.155605333333333
    "H17" = 3.466816
    "I17" = 0.04233532586138062
    "J17" = 0.04233532586138063
    "M17" = 0.5052413333333333
    "N17" = 1.515724
    "O17" = 0.01234989497255689
    "P17" = 0.01234989497255689
    "Q17" = 0.5838595794204443
    "R17" = 5.254736214784
    "S17" = 0.0005228368280170223
    "T17" = 0.0005228368280170224
    "G18" = 1.155605333333333
    "H18" = 3.466816
    "I18" = 0.04233532586138062
    "J18" = 0.04233532586138063
    "M18" = 11.46238066666667
    "N18" = 34.387142
    "O18" = 0.280181347070047
    "P18" = 0.280181347070047
    "Q18" = 13.24598823109689
    "R18" = 119.213894079872
    "S18" = 0.01186156862849102
    "T18" = 0.01186156862849102
    "G19" = 1.155605333333333
    "H19" = 3.466816
    "I19" = 0.04233532586138062
    "J19" = 0.04233532586138063
    "M19" = 5.621112333333333
    "N19" = 16.863337
    "O19" = 0.1373999757454739
    "P19" = 0.1373999757454739
    "Q19" = 6.495787391665777
    "R19" = 58.462086524992
    "S19" = 0.005816872746530431
    "T19" = 0.005816872746530433
}

foreach ($key in $updates.Keys) {
    $ws.Range($key).Value = $updates[$key]
}
